$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q0) updates
$ws.Range("B3").Value = 0.1068730318990097
$ws.Range("C3").Value = 0.6104158195844518
$ws.Range("D3").Value = 0.6472573004727202
$ws.Range("E3").Value = 0.8045230266889322
$ws.Range("F3").Value = 0.8253800783743511
$ws.Range("G3").Value = 15

# Row 4 (Q1) updates
$ws.Range("B4").Value = 0.2448706236272014
$ws.Range("C4").Value = 0.7412646195823523
$ws.Range("D4").Value = 0.8458750405844518
$ws.Range("E4").Value = 0.9197146517178314
$ws.Range("F4").Value = 0.9199827697284622
$ws.Range("G4").Value = 14
